$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.041.09'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.47%  '
$ws.Range('E2').Style = 'Normal'

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.829.95'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.34%  '
$ws.Range('E3').Style = 'Normal'

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9987'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('E4').Style = 'Normal'

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.27'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.33%  '
$ws.Range('E5').Style = 'Normal'

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6283'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -4.91%  '
$ws.Range('E6').Style = 'Normal'

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.0000'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E7').Style = 'Normal'

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07614'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +2.52%  '
$ws.Range('E8').Style = 'Normal'

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.17'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +7.75%  '
$ws.Range('E9').Style = 'Normal'

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.2913'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.81%  '
$ws.Range('E10').Style = 'Normal'

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '22.78'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.65%  '
$ws.Range('E11').Style = 'Normal'

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07648'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.45%  '
$ws.Range('E12').Style = 'Normal'

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.829.85'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.12%  '
$ws.Range('E13').Style = 'Normal'

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.960'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.74%  '
$ws.Range('E14').Style = 'Normal'

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6648'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.38%  '
$ws.Range('E15').Style = 'Normal'

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '82.47'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.57%  '
$ws.Range('E16').Style = 'Normal'

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000009369'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +9.45%  '
$ws.Range('E17').Style = 'Normal'

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.990'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.99%  '
$ws.Range('E18').Style = 'Normal'

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '28.857.15'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.08%  '
$ws.Range('E19').Style = 'Normal'

# Row 20
$ws.Range('B20').NumberFormat = '@'
$ws.Range('B20').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('B20').Style = 'Normal'
$ws.Range('C20').NumberFormat = '@'
$ws.Range('C20').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('C20').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.072.00'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.82%  '
$ws.Range('E20').Style = 'Normal'

# Row 21
$ws.Range('B21').NumberFormat = '@'
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('B21').Style = 'Normal'
$ws.Range('C21').NumberFormat = '@'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('C21').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '225.17'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.75%  '
$ws.Range('E21').Style = 'Normal'

# Row 22
$ws.Range('B22').NumberFormat = '@'
$ws.Range('B22').Value = 'Avalanche'
$ws.Range('B22').Style = 'Normal'
$ws.Range('C22').NumberFormat = '@'
$ws.Range('C22').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('C22').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.32'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.18%  '
$ws.Range('E22').Style = 'Normal'

# Row 23
$ws.Range('B23').NumberFormat = '@'
$ws.Range('B23').Value = 'Dai'
$ws.Range('B23').Style = 'Normal'
$ws.Range('C23').NumberFormat = '@'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('C23').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9992'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('E23').Style = 'Normal'

# Row 24
$ws.Range('B24').NumberFormat = '@'
$ws.Range('B24').Value = 'Chainlink'
$ws.Range('B24').Style = 'Normal'
$ws.Range('C24').NumberFormat = '@'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('C24').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.218'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.72%  '
$ws.Range('E24').Style = 'Normal'

# Row 25
$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value = 'BinanceUSD'
$ws.Range('B25').Style = 'Normal'
$ws.Range('C25').NumberFormat = '@'
$ws.Range('C25').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('C25').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.000'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('E25').Style = 'Normal'

# Row 26
$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'Monero'
$ws.Range('B26').Style = 'Normal'
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('C26').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '159.91'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('E26').Style = 'Normal'

# Row 27
$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('B27').Style = 'Normal'
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('C27').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.423'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.36%  '
$ws.Range('E27').Style = 'Normal'

# Row 28
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'Stellar'
$ws.Range('B28').Style = 'Normal'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('C28').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1365'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.53%  '
$ws.Range('E28').Style = 'Normal'

# Row 29
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('B29').Style = 'Normal'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('C29').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '17.82'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.83%  '
$ws.Range('E29').Style = 'Normal'

# Row 30
$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('B30').Style = 'Normal'
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('C30').Style = 'Normal'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.493'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.71%  '
$ws.Range('E30').Style = 'Normal'

# Row 32
$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('B32').Style = 'Normal'
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('C32').Style = 'Normal'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.031'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.31%  '
$ws.Range('E32').Style = 'Normal'

# Row 33
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'Toncoin'
$ws.Range('B33').Style = 'Normal'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('C33').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.202'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.91%  '
$ws.Range('E33').Style = 'Normal'

# Row 34
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'Hedera'
$ws.Range('B34').Style = 'Normal'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('C34').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05206'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.31%  '
$ws.Range('E34').Style = 'Normal'

# Row 35
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('B35').Style = 'Normal'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('C35').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.842'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.21%  '
$ws.Range('E35').Style = 'Normal'

# Row 36
$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('B36').Style = 'Normal'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('C36').Style = 'Normal'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.154'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.71%  '
$ws.Range('E36').Style = 'Normal'

# Row 37
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('B37').Style = 'Normal'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('C37').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.7315'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.51%  '
$ws.Range('E37').Style = 'Normal'

# Row 38
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'HuobiToken'
$ws.Range('B38').Style = 'Normal'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('C38').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.612'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.69%  '
$ws.Range('E38').Style = 'Normal'

# Row 39
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'Maker'
$ws.Range('B39').Style = 'Normal'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('C39').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.275.47'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.90%  '
$ws.Range('E39').Style = 'Normal'

# Row 40
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'MXToken'
$ws.Range('B40').Style = 'Normal'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('C40').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.760'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.88%  '
$ws.Range('E40').Style = 'Normal'

# Row 41
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'VeChain'
$ws.Range('B41').Style = 'Normal'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('C41').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01786'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.41%  '
$ws.Range('E41').Style = 'Normal'

# Row 42
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('B42').Style = 'Normal'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('C42').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.500'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +7.35%  '
$ws.Range('E42').Style = 'Normal'

# Row 43
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8884'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -3.38%  '
$ws.Range('E43').Style = 'Normal'

# Row 44
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('B44').Style = 'Normal'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('C44').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.000'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.12%  '
$ws.Range('E44').Style = 'Normal'

# Row 45
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'Quant'
$ws.Range('B45').Style = 'Normal'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('C45').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '101.45'
$ws.Range('D45').Style = 'Normal'

# Row 46
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'RocketPoolETH'
$ws.Range('B46').Style = 'Normal'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('C46').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.972.77'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.09%  '
$ws.Range('E46').Style = 'Normal'

# Row 47
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'Mantle'
$ws.Range('B47').Style = 'Normal'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('C47').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5105'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.76%  '
$ws.Range('E47').Style = 'Normal'

# Row 48
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'Aave'
$ws.Range('B48').Style = 'Normal'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('C48').Style = 'Normal'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '63.74'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.19%  '
$ws.Range('E48').Style = 'Normal'

# Row 49
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'XinFinNetwork'
$ws.Range('B49').Style = 'Normal'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range('C49').Style = 'Normal'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.07417'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -12.68%  '
$ws.Range('E49').Style = 'Normal'

# Row 50
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'TheSandbox'
$ws.Range('B50').Style = 'Normal'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('C50').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.3980'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.75%  '
$ws.Range('E50').Style = 'Normal'

# Row 51
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('B51').Style = 'Normal'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('C51').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.827'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.11%  '
$ws.Range('E51').Style = 'Normal'
